# Update styling and key-result page
# -----------------------------------
# 1) Correct the reportingManagerId (column H) values for a handful of
#    employee rows.
# 2) Append two new rows (21-22) to the "key result" / id-lookup tail of
#    the sheet, following the same blank-template pattern already used by
#    rows 19-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix reportingManagerId values -------------------------------------
$ws.Range("H6").Value  = 2
$ws.Range("H9").Value  = 2
$ws.Range("H10").Value = 5
$ws.Range("H11").Value = 5
$ws.Range("H12").Value = 8
$ws.Range("H13").Value = 8

# --- 2. Append row 21 -------------------------------------------------------
# Clone the formatting of the existing "blank template" row (20) so every
# column keeps a real (but empty) cell, then fill in the values that differ.
$ws.Range("A20:K20").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A21").Formula = "="""""         # empty-text placeholder
$ws.Range("B21").Formula = "="""""
$ws.Range("C21").Value   = "[]"
$ws.Range("E21").Value   = "[]"
$ws.Range("G21").Formula = "="""""
$ws.Range("I21").Value   = 20
$ws.Range("J21").Formula = "="""""
$ws.Range("K21").Formula = "="""""

# --- 3. Append row 22 -------------------------------------------------------
$ws.Range("A20:K20").Copy()
$ws.Range("A22:K22").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C22").Value = "[]"
$ws.Range("E22").Value = "[]"
$ws.Range("I22").Value = 21

$excel.CutCopyMode = $false
